$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several candidates' "Date Dropped" cell (column G) currently holds the
# placeholder text "-". Update them to the actual dropout dates, recorded
# as real dates (matching the existing date-formatted cells in the column).
#
# Copy the date format already used by the other "Date Dropped" cells (e.g.
# G3) onto the placeholder cells first, so they pick up the existing date
# style instead of Excel minting a brand new one, then set the values.
$ws.Range("G3").Copy()
$ws.Range("G2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G8").PasteSpecial(-4122)
$ws.Range("G17").PasteSpecial(-4122)
$ws.Range("G23").PasteSpecial(-4122)

$ws.Range("G2").Value = 43872    # 2/11/2020 - Bennet
$ws.Range("G8").Value = 43891    # 3/1/2020  - Gabbard
$ws.Range("G17").Value = 43892   # 3/2/2020  - Ryan
$ws.Range("G23").Value = 43872   # 2/11/2020 - Hickenlooper

# Reflect where the editor left the selection.
$ws.Range("I18").Select()
